# Update cryptocurrency price/volume figures (D: Price, E: Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.167.87"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "2.550.88"
$ws.Range("E3").Value = "  -2.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.22%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").Value = "2.549.71"
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("E11").Value = "  +1.84%  "

$ws.Range("E12").Value = "  -0.77%  "

$ws.Range("E13").Value = "  -4.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "3.017.42"
$ws.Range("E15").Value = "  -2.43%  "

$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").Value = "67.060.81"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").Value = "2.542.39"
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.66"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("E24").Value = "  +5.90%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  +0.90%  "

$ws.Range("E27").Value = "  -3.86%  "

$ws.Range("D28").Value = "2.673.07"
$ws.Range("E28").Value = "  -3.08%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "534.58"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.53%  "

$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.09"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("E39").Value = "  -0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.46"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.357"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.70%  "

$ws.Range("E42").Value = "  +0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.99%  "

$ws.Range("E44").Value = "  +5.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.563"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.35%  "

$ws.Range("D49").Value = "0.0₆0279"
$ws.Range("E49").Value = "  -5.24%  "

$ws.Range("E50").Value = "  -1.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.19%  "
